# Weekly refresh of the Perejil (Vega Monumental Concepción) price series.
# A new week's observation (pair of rows: "Primera"/"Segunda" quality) is
# published, so the whole data block shifts down by one pair of rows and the
# newest pair is inserted right after the header/previous block, reusing the
# most recent price pattern as the starting point for the new date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data (rows 142:189) down by two rows to make room for
# the new weekly observation.
$ws.Rows.Item(142).Insert()
$ws.Rows.Item(142).Insert()

# After the insert, the last observation that used to sit at 188:189 now
# lives at 190:191. Duplicate it into the freshly opened 142:143 slot so the
# new entry starts from the latest known price/unit/origin pattern.
$ws.Range("A190:R191").Copy()
$ws.Range("A142").PasteSpecial()
$excel.CutCopyMode = $false

# Stamp the new pair with the new reporting date.
$ws.Range("D142").Value = 44924
$ws.Range("D143").Value = 44924
